$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.964.74'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.26'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.76'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5065'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3660'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07202'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8942'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.72'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.84'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07528'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '95.05'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.241'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008546'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.25'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9997'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.016.21'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.033'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.114.81'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.41'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.420'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.46'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.783'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.19%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.079'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.45'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.707'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.687'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09160'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05149'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7533'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.989'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.160'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.228'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.571'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5658'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.18%  '
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.071'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.612'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.75'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.536'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1477'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4737'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.0000'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.12'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.566'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.89'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.22'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.00%  '
